$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 588.7857
$ws.Range("I33").Value = 248.2
$ws.Range("K33").Value = 248.2
$ws.Range("M33").Value = -19.19999999999999
$ws.Range("H43").Value = 5613.357
$ws.Range("I43").Value = 5658.6
$ws.Range("J43").Value = 5588.222
$ws.Range("K43").Value = 5658.6
$ws.Range("L43").Value = 5588.222
$ws.Range("M43").Value = -5589.6
$ws.Range("N43").Value = -5726.222
$ws.Range("H70").Value = 127190.75
$ws.Range("J70").Value = 202116.2
$ws.Range("L70").Value = 606348.6000000001
$ws.Range("N70").Value = -606888.6000000001
$ws.Range("H73").Value = 127190.75
$ws.Range("J73").Value = 202116.2
$ws.Range("L73").Value = 606348.6000000001
$ws.Range("N73").Value = -608220.6000000001
$ws.Range("H106").Value = 1616.2
$ws.Range("I106").Value = 1616.2
$ws.Range("K106").Value = 1616.2
$ws.Range("M106").Value = -985.2
$ws.Range("H136").Value = 65997.664
$ws.Range("J136").Value = 65997.664
$ws.Range("L136").Value = 65997.664
$ws.Range("N136").Value = -76197.664
$ws.Range("H137").Value = 4999
$ws.Range("J137").Value = 4998.5
$ws.Range("L137").Value = 14995.5
$ws.Range("N137").Value = -20095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7605.6562
$ws.Range("I61").Value = 5145.84
$ws.Range("J61").Value = 16390.715
$ws.Range("K61").Value = 5145.84
$ws.Range("L61").Value = 16390.715
$ws.Range("M61").Value = -4933.84
$ws.Range("N61").Value = -16814.715
$ws.Range("H74").Value = 1861.909
$ws.Range("I74").Value = 1650.2941
$ws.Range("J74").Value = 2086.75
$ws.Range("K74").Value = 1650.2941
$ws.Range("L74").Value = 2086.75
$ws.Range("M74").Value = -776.2941000000001
$ws.Range("N74").Value = -3834.75
$ws.Range("H77").Value = 1861.909
$ws.Range("I77").Value = 1650.2941
$ws.Range("J77").Value = 2086.75
$ws.Range("K77").Value = 8251.470499999999
$ws.Range("L77").Value = 10433.75
$ws.Range("M77").Value = -3883.470499999999
$ws.Range("N77").Value = -19169.75
$ws.Range("H122").Value = 3552.3333
$ws.Range("I122").Value = 3360.3333
$ws.Range("K122").Value = 10080.9999
$ws.Range("M122").Value = -7630.999899999999
$ws.Range("H136").Value = 7605.6562
$ws.Range("I136").Value = 5145.84
$ws.Range("J136").Value = 16390.715
$ws.Range("K136").Value = 15437.52
$ws.Range("L136").Value = 49172.145
$ws.Range("M136").Value = -12887.52
$ws.Range("N136").Value = -54272.145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2567.3333
$ws.Range("I20").Value = 2766.1538
$ws.Range("K20").Value = 2766.1538
$ws.Range("M20").Value = -2519.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9670.299999999999
$ws.Range("I31").Value = 20225
$ws.Range("K31").Value = 20225
$ws.Range("M31").Value = -19930
$ws.Range("H34").Value = 9670.299999999999
$ws.Range("I34").Value = 20225
$ws.Range("K34").Value = 20225
$ws.Range("M34").Value = -20023
$ws.Range("H58").Value = 3701
$ws.Range("I58").Value = 3101.375
$ws.Range("J58").Value = 6099.5
$ws.Range("K58").Value = 3101.375
$ws.Range("L58").Value = 6099.5
$ws.Range("M58").Value = -2898.375
$ws.Range("N58").Value = -6505.5
$ws.Range("H62").Value = 10886.2
$ws.Range("J62").Value = 13053.25
$ws.Range("L62").Value = 13053.25
$ws.Range("N62").Value = -14301.25
$ws.Range("H65").Value = 10886.2
$ws.Range("J65").Value = 13053.25
$ws.Range("L65").Value = 65266.25
$ws.Range("N65").Value = -71506.25
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H81").Value = 38327.055
$ws.Range("J81").Value = 38327.055
$ws.Range("L81").Value = 38327.055
$ws.Range("N81").Value = -40323.055
$ws.Range("H84").Value = 38327.055
$ws.Range("J84").Value = 38327.055
$ws.Range("L84").Value = 114981.165
$ws.Range("N84").Value = -124965.165
$ws.Range("H132").Value = 4884.2812
$ws.Range("I132").Value = 3742.5
$ws.Range("J132").Value = 8309.625
$ws.Range("K132").Value = 11227.5
$ws.Range("L132").Value = 24928.875
$ws.Range("M132").Value = -8697.5
$ws.Range("N132").Value = -29988.875
$ws.Range("H134").Value = 4426.0386
$ws.Range("I134").Value = 3958.7727
$ws.Range("K134").Value = 11876.3181
$ws.Range("M134").Value = -9341.3181
$ws.Range("H136").Value = 3701
$ws.Range("I136").Value = 3101.375
$ws.Range("J136").Value = 6099.5
$ws.Range("K136").Value = 9304.125
$ws.Range("L136").Value = 18298.5
$ws.Range("M136").Value = -6754.125
$ws.Range("N136").Value = -23398.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 699.5
$ws.Range("I68").Value = 699
$ws.Range("J68").Value = 700
$ws.Range("K68").Value = 2097
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -1286
$ws.Range("N68").Value = -3722
$ws.Range("H71").Value = 699.5
$ws.Range("I71").Value = 699
$ws.Range("J71").Value = 700
$ws.Range("K71").Value = 6291
$ws.Range("L71").Value = 700
$ws.Range("M71").Value = -2235
$ws.Range("N71").Value = -14412
$ws.Range("H98").Value = 1295.0714
$ws.Range("J98").Value = 1589.6666
$ws.Range("L98").Value = 4768.9998
$ws.Range("N98").Value = -7764.9998
$ws.Range("H107").Value = 1177.2222
$ws.Range("J107").Value = 1999.25
$ws.Range("L107").Value = 5997.75
$ws.Range("N107").Value = -9837.75
$ws.Range("H122").Value = 33999.332
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 33999.332
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 305993.988
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -310893.988

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7466.1333
$ws.Range("I80").Value = 6982.4287
$ws.Range("K80").Value = 6982.4287
$ws.Range("M80").Value = -5984.4287
$ws.Range("H83").Value = 7466.1333
$ws.Range("I83").Value = 6982.4287
$ws.Range("K83").Value = 34912.14350000001
$ws.Range("M83").Value = -29920.14350000001
$ws.Range("H113").Value = 38384.855
$ws.Range("I113").Value = 32499.75
$ws.Range("J113").Value = 46231.668
$ws.Range("K113").Value = 32499.75
$ws.Range("L113").Value = 46231.668
$ws.Range("M113").Value = -30329.75
$ws.Range("N113").Value = -50571.668
$ws.Range("H122").Value = 4602.375
$ws.Range("I122").Value = 3138.6
$ws.Range("K122").Value = 9415.799999999999
$ws.Range("M122").Value = -6965.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2164.7097
$ws.Range("I22").Value = 2083.3572
$ws.Range("J22").Value = 2231.7058
$ws.Range("K22").Value = 2083.3572
$ws.Range("L22").Value = 2231.7058
$ws.Range("M22").Value = -1788.3572
$ws.Range("N22").Value = -2821.7058
$ws.Range("H27").Value = 2164.7097
$ws.Range("I27").Value = 2083.3572
$ws.Range("J27").Value = 2231.7058
$ws.Range("K27").Value = 2083.3572
$ws.Range("L27").Value = 2231.7058
$ws.Range("M27").Value = -1976.3572
$ws.Range("N27").Value = -2445.7058
$ws.Range("H46").Value = 4319.4
$ws.Range("I46").Value = 3932.6667
$ws.Range("K46").Value = 3932.6667
$ws.Range("M46").Value = -3744.6667
$ws.Range("H68").Value = 6365
$ws.Range("I68").Value = 6925.8335
$ws.Range("K68").Value = 6925.8335
$ws.Range("M68").Value = -6176.8335
$ws.Range("H71").Value = 6365
$ws.Range("I71").Value = 6925.8335
$ws.Range("K71").Value = 34629.1675
$ws.Range("M71").Value = -30885.1675
$ws.Range("H82").Value = 1608.9231
$ws.Range("I82").Value = 1545.3
$ws.Range("J82").Value = 1821
$ws.Range("K82").Value = 1545.3
$ws.Range("L82").Value = 1821
$ws.Range("M82").Value = -1184.3
$ws.Range("N82").Value = -2543
$ws.Range("H85").Value = 1608.9231
$ws.Range("I85").Value = 1545.3
$ws.Range("J85").Value = 1821
$ws.Range("K85").Value = 1545.3
$ws.Range("L85").Value = 1821
$ws.Range("M85").Value = -297.3
$ws.Range("N85").Value = -4317
$ws.Range("H136").Value = 1442.3529
$ws.Range("I136").Value = 1276.6666
$ws.Range("K136").Value = 3829.9998
$ws.Range("M136").Value = -1279.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 15000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 15000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -15284
$ws.Range("H17").Value = 7500
$ws.Range("J17").Value = 7500
$ws.Range("L17").Value = 7500
$ws.Range("N17").Value = -7844
$ws.Range("H122").Value = 4191.1113
$ws.Range("I122").Value = 2953.4167
$ws.Range("K122").Value = 8860.250100000001
$ws.Range("M122").Value = -6410.250100000001
$ws.Range("H132").Value = 2434.2273
$ws.Range("I132").Value = 2463.85
$ws.Range("K132").Value = 7391.549999999999
$ws.Range("M132").Value = -4861.549999999999
$ws.Range("H136").Value = 30000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 30000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 90000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -95100
$ws.Range("H138").Value = 108333
$ws.Range("J138").Value = 117499.5
$ws.Range("L138").Value = 117499.5
$ws.Range("N138").Value = -127779.5
$ws.Range("H139").Value = 78713
$ws.Range("I139").Value = 78713
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 78713
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -73573
$ws.Range("N139").ClearContents()
